$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# 1. Title
$d.Content.Find.Execute("The Elusive Beauty of Black Holes", $false, $true, $false, $false, $false, $true, 1, $false, "Unveiling the Symphony of Chemistry: Fundamental Principles and Practical Applications", 2) | Out-Null

# 2. Author name
$d.Content.Find.Execute("Lilith Sinclair", $false, $true, $false, $false, $false, $true, 1, $false, "Alma Einstein", 2) | Out-Null

# 3. username (match whole word to avoid touching 'Lilith' again, case-sensitive lowercase)
$d.Content.Find.Execute("lilith", $true, $true, $false, $false, $false, $true, 1, $false, "alma", 2) | Out-Null

# 4. email domain part
$d.Content.Find.Execute("sinclair@celestialscience", $false, $true, $false, $false, $false, $true, 1, $false, "einstein@eduverse", 2) | Out-Null

# 5. Big paragraph body (3 sections)
$old5 = 'In the enigmatic depths of the cosmos, where light surrenders to gravitational allure, lies an enigmatic entity - the black hole. These celestial maelstroms, born from the cataclysmic demise of massive stars, hold the power to warp the fabric of space-time, twisting and distorting the surrounding universe. With their event horizons serving as cosmic borders, they become realms of mystery, gateways to a world of phenomena that challenge our current scientific understanding. Black holes allure and mystify, beckoning us to unravel their profound secrets and explore their tantalizing enigmas.As we delve into the captivating realm of black hole physics, we confront an array of perplexing questions that ignite our curiosity. What lies within these cosmic abysses? How can we penetrate their event horizons and gain insights into the enigmatic physics that govern their behavior? Do black holes serve as cosmic gateways to other dimensions, unlocking avenues of travel across vast interstellar distances? These inquiries push the boundaries of human knowledge, propelling us to explore the cosmos in unprecedented ways.Black holes hold a mirror to the very foundation of physics, unveiling the interplay between gravity, quantum mechanics, and thermodynamics. Their enigmatic nature confronts our current scientific paradigms, compelling us to re-examine deeply held beliefs and embark on an audacious voyage of discovery. By unveiling the secrets of black holes, we illuminate the path towards a deeper understanding of the universe and secure our place in the cosmos as sentient beings, striving to unravel the vast tapestry of reality.'
$new5 = 'Journey through Chemistry''s Captivating Realm:^lChemistry, a captivating branch of science, offers a profound understanding of the material world, unlocking the innermost secrets of matter and its interactions. It delves into the realm of atoms and molecules, unraveling the fundamental principles that govern their behavior and shape the world around us. Chemistry has far-reaching applications, impacting various aspects of our lives, from the creation of life-saving medicines and innovative technologies to the development of sustainable energy sources.^l^lChemistry''s Role in Shaping Our World:^lThe field of chemistry plays a pivotal role in numerous industries, contributing to the progress of society. It drives advancements in medicine by fostering the discovery of new drugs and treatments, paving the way for a healthier future. Chemistry also fuels innovation in technology, enabling the development of advanced materials, energy-efficient processes, and eco-friendly products, propelling us toward a more sustainable existence.^l^lExploring Chemistry''s Fascinating Phenomena:^lChemistry captivates us with its intriguing phenomena, revealing the hidden forces at play in the world around us. From the mesmerizing reactions that produce vibrant colors and produce intense heat to the complex interactions that result in the synthesis of new substances, chemistry offers endless opportunities for exploration and discovery. These phenomena inspire a sense of wonder and curiosity, igniting a passion for understanding the complexities of the natural world.'
$d.Content.Find.Execute($old5, $false, $true, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# 6. Summary paragraph body
$old7 = 'The profound allure of black holes stems from their ability to challenge our current scientific knowledge and ignite our insatiable curiosity. These mysterious entities serve as cosmic laboratories where fundamental physics is put to the test, compelling us to delve deeper into the intricate tapestry of the universe. Through the study of black holes, we glimpse the profound interplay between gravity, quantum mechanics, and thermodynamics, pushing the boundaries of our scientific understanding and expanding our perception of the cosmos.'
$new7 = 'Chemistry is a captivating field of science that unveils the symphony of matter and its interactions. It plays a crucial role in shaping our world through its far-reaching applications in medicine, technology, and sustainability. Chemistry captivates us with its fascinating phenomena, stimulating curiosity and igniting a passion for understanding the intricacies of the natural world. Exploring chemistry opens doors to a world of discovery and understanding, empowering us to solve complex challenges and shape a better future.'
$d.Content.Find.Execute($old7, $false, $true, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# 7. Add a trailing empty paragraph at the end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
